$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6
$ws.Range("E6").Value = ""
$ws.Range("H6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 8
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H8").Value = "15,23 TL - 30,47 TL - 304,72 TL"

# Row 9
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H9").Value = "15,23 TL - 30,47 TL - 304,72 TL"

# Row 10
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H10").Value = "15,23 TL - 30,47 TL - 304,72 TL"

# Row 11
$ws.Range("E11").Value = ""
$ws.Range("H11").Value = "3,05 TL - 6,1 TL - 76,18 TL"

# Row 13
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"
$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
$ws.Range("H14").Value = "2.100 TL - 4.300 TL"

# Row 15
$ws.Range("D15").Value = "%0,8 Asgari Tutar: 285,72 TL Azami Tutar: 285,72 TL / 3.047,62 TL"

# Row 17
$ws.Range("D17").Value = "%0,8 Asgari Tutar: 285,72 TL Azami Tutar: 285,72 TL / 3.047,62 TL"

# Row 20
$ws.Range("D20").Value = "285,72 TL"

# Row 21
$ws.Range("D21").Value = "%0,5 Asgari Tutar: 428,58 TL Azami Tutar: 428,58 TL / 5.523,81 TL"

# Row 22
$ws.Range("D22").Value = "%1 Asgari Tutar: 285,72 TL Azami Tutar: 285,72 TL / 6.857,15 TL"

# Row 23
$ws.Range("D23").Value = "64,77 TL"

# Row 24
$ws.Range("D24").Value = "476,2 TL"

# Row 25
$ws.Range("D25").Value = "428,58 TL"
